$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New stimuli order (image / word / category) for rows 2-33,
# one entry per row, in row order.
$stimuli = @(
    ,@("dog/dog106.png", "spenden", "dog")
    ,@("dog/dog101.png", "bauen", "dog")
    ,@("flower/flower076.png", "angeln", "flower")
    ,@("dog/dog086.png", "lügen", "dog")
    ,@("flower/flower091.png", "lernen", "flower")
    ,@("flower/flower080.png", "lächeln", "flower")
    ,@("flower/flower092.png", "kennen", "flower")
    ,@("flower/flower097.png", "proben", "flower")
    ,@("flower/flower071.png", "mögen", "flower")
    ,@("dog/dog080.png", "heben", "dog")
    ,@("flower/flower074.png", "hören", "flower")
    ,@("dog/dog087.png", "holen", "dog")
    ,@("dog/dog078.png", "kriegen", "dog")
    ,@("flower/flower070.png", "ärgern", "flower")
    ,@("dog/dog092.png", "trotzen", "dog")
    ,@("dog/dog108.png", "münzen", "dog")
    ,@("dog/dog073.png", "spüren", "dog")
    ,@("dog/dog117.png", "binden", "dog")
    ,@("dog/dog072.png", "narren", "dog")
    ,@("dog/dog104.png", "quellen", "dog")
    ,@("dog/dog094.png", "stoppen", "dog")
    ,@("flower/flower096.png", "parken", "flower")
    ,@("flower/flower107.png", "prüfen", "flower")
    ,@("dog/dog065.png", "zielen", "dog")
    ,@("flower/flower084.png", "grenzen", "flower")
    ,@("flower/flower068.png", "kranken", "flower")
    ,@("dog/dog084.png", "lassen", "dog")
    ,@("dog/dog069.png", "heißen", "dog")
    ,@("flower/flower110.png", "rufen", "flower")
    ,@("flower/flower109.png", "achten", "flower")
    ,@("flower/flower104.png", "betteln", "flower")
    ,@("flower/flower108.png", "wachsen", "flower")
)

for ($i = 0; $i -lt $stimuli.Length; $i++) {
    $row = $i + 2
    $entry = $stimuli[$i]
    $ws.Cells.Item($row, 2).Value = $entry[0]   # B: image
    $ws.Cells.Item($row, 3).Value = $entry[1]   # C: word
    $ws.Cells.Item($row, 4).Value = $entry[2]   # D: category
}
